$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H19").Value = 670.2857
$ws.Range("I19").Value = 467
$ws.Range("J19").Value = 822.75
$ws.Range("K19").Value = 467
$ws.Range("L19").Value = 822.75
$ws.Range("M19").Value = -292
$ws.Range("N19").Value = -1172.75
$ws.Range("H98").Value = 823.1070999999999
$ws.Range("I98").Value = 854.43475
$ws.Range("J98").Value = 679
$ws.Range("K98").Value = 854.43475
$ws.Range("L98").Value = 679
$ws.Range("M98").Value = 643.56525
$ws.Range("N98").Value = -3675
$ws.Range("H122").Value = 823.1070999999999
$ws.Range("I122").Value = 854.43475
$ws.Range("J122").Value = 679
$ws.Range("K122").Value = 2563.30425
$ws.Range("L122").Value = 2037
$ws.Range("M122").Value = -113.3042500000001
$ws.Range("N122").Value = -6937
$ws.Range("H129").Value = 884.86
$ws.Range("I129").Value = 478.77777
$ws.Range("J129").Value = 974
$ws.Range("K129").Value = 1436.33331
$ws.Range("L129").Value = 2922
$ws.Range("M129").Value = 3563.66669
$ws.Range("N129").Value = -12922
$ws.Range("H135").Value = 2112.1875
$ws.Range("I135").Value = 1846.72
$ws.Range("J135").Value = 3060.2856
$ws.Range("K135").Value = 16620.48
$ws.Range("L135").Value = 27542.5704
$ws.Range("M135").Value = -14085.48
$ws.Range("N135").Value = -32612.5704
$ws.Range("H137").Value = 1163.725
$ws.Range("I137").Value = 1119.5161
$ws.Range("J137").Value = 1316
$ws.Range("K137").Value = 3358.5483
$ws.Range("L137").Value = 3948
$ws.Range("M137").Value = -808.5483000000004
$ws.Range("N137").Value = -9048
$ws.Range("H141").Value = 1735
$ws.Range("I141").Value = 1618.3928
$ws.Range("J141").Value = 2551.25
$ws.Range("K141").Value = 4855.178400000001
$ws.Range("L141").Value = 7653.75
$ws.Range("M141").Value = 324.8215999999993
$ws.Range("N141").Value = -18013.75

$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 1448.94
$ws.Range("I61").Value = 1138.641
$ws.Range("J61").Value = 2549.0908
$ws.Range("K61").Value = 1138.641
$ws.Range("L61").Value = 2549.0908
$ws.Range("M61").Value = -926.6410000000001
$ws.Range("N61").Value = -2973.0908
$ws.Range("H74").Value = 21740736
$ws.Range("I74").Value = 1601.2
$ws.Range("J74").Value = 62501612
$ws.Range("K74").Value = 1601.2
$ws.Range("L74").Value = 62501612
$ws.Range("M74").Value = -727.2
$ws.Range("N74").Value = -62503360
$ws.Range("H77").Value = 21740736
$ws.Range("I77").Value = 1601.2
$ws.Range("J77").Value = 62501612
$ws.Range("K77").Value = 8006
$ws.Range("L77").Value = 312508060
$ws.Range("M77").Value = -3638
$ws.Range("N77").Value = -312516796
$ws.Range("H132").Value = 2856.75
$ws.Range("I132").Value = 1794.963
$ws.Range("K132").Value = 5384.889
$ws.Range("M132").Value = -2854.889
$ws.Range("H134").Value = 44420
$ws.Range("J134").Value = 44420
$ws.Range("L134").Value = 44420
$ws.Range("N134").Value = -54560
$ws.Range("H136").Value = 1448.94
$ws.Range("I136").Value = 1138.641
$ws.Range("J136").Value = 2549.0908
$ws.Range("K136").Value = 3415.923
$ws.Range("L136").Value = 7647.2724
$ws.Range("M136").Value = -865.9230000000002
$ws.Range("N136").Value = -12747.2724
$ws.Range("H140").Value = 40429
$ws.Range("J140").Value = 40429
$ws.Range("L140").Value = 40429
$ws.Range("N140").Value = -50789

$ws = $wb.Worksheets("BSM")
$ws.Range("H134").Value = 1578.0428
$ws.Range("I134").Value = 1355.7059
$ws.Range("J134").Value = 2174.842
$ws.Range("K134").Value = 4067.1177
$ws.Range("L134").Value = 6524.526
$ws.Range("M134").Value = -1532.1177
$ws.Range("N134").Value = -11594.526

$ws = $wb.Worksheets("CRP")
$ws.Range("H105").Value = 14494660
$ws.Range("I105").Value = 20835504
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 20835504
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = -20833757
$ws.Range("N105").Value = -4794

$ws = $wb.Worksheets("CUL")
$ws.Range("H123").Value = 5162.727
$ws.Range("I123").Value = 895
$ws.Range("J123").Value = 6111.1113
$ws.Range("K123").Value = 2685
$ws.Range("L123").Value = 18333.3339
$ws.Range("M123").Value = -235
$ws.Range("N123").Value = -23233.3339

$ws = $wb.Worksheets("GSM")
$ws.Range("H11").Value = 295857150
$ws.Range("I11").Value = 504000000
$ws.Range("J11").Value = 18333334
$ws.Range("K11").Value = 504000000
$ws.Range("L11").Value = 18333334
$ws.Range("M11").Value = -503999861
$ws.Range("N11").Value = -18333612
$ws.Range("H80").Value = 2740.6897
$ws.Range("I80").Value = 2498.5715
$ws.Range("J80").Value = 2966.6667
$ws.Range("K80").Value = 2498.5715
$ws.Range("L80").Value = 2966.6667
$ws.Range("M80").Value = -1500.5715
$ws.Range("N80").Value = -4962.6667
$ws.Range("H83").Value = 2740.6897
$ws.Range("I83").Value = 2498.5715
$ws.Range("J83").Value = 2966.6667
$ws.Range("K83").Value = 12492.8575
$ws.Range("L83").Value = 14833.3335
$ws.Range("M83").Value = -7500.8575
$ws.Range("N83").Value = -24817.3335
$ws.Range("H132").Value = 6413603.5
$ws.Range("I132").Value = 9807414
$ws.Range("K132").Value = 29422242
$ws.Range("M132").Value = -29419712

$ws = $wb.Worksheets("LTW")
$ws.Range("H9").Value = 1000
$ws.Range("J9").Value = 1000
$ws.Range("L9").Value = 1000
$ws.Range("N9").Value = -1448
$ws.Range("H13").Value = 7210
$ws.Range("J13").Value = 7210
$ws.Range("L13").Value = 7210
$ws.Range("N13").Value = -7490
$ws.Range("H82").Value = 1228688
$ws.Range("I82").Value = 3334666.2
$ws.Range("J82").Value = 175698.83
$ws.Range("K82").Value = 3334666.2
$ws.Range("L82").Value = 175698.83
$ws.Range("M82").Value = -3334305.2
$ws.Range("N82").Value = -176420.83
$ws.Range("H85").Value = 1228688
$ws.Range("I85").Value = 3334666.2
$ws.Range("J85").Value = 175698.83
$ws.Range("K85").Value = 3334666.2
$ws.Range("L85").Value = 175698.83
$ws.Range("M85").Value = -3333418.2
$ws.Range("N85").Value = -178194.83
$ws.Range("H132").Value = 9168908
$ws.Range("I132").Value = 13482511
$ws.Range("J132").Value = 2501.875
$ws.Range("K132").Value = 40447533
$ws.Range("L132").Value = 7505.625
$ws.Range("M132").Value = -40445003
$ws.Range("N132").Value = -12565.625
$ws.Range("H134").Value = 52420
$ws.Range("J134").Value = 52420
$ws.Range("L134").Value = 52420
$ws.Range("N134").Value = -62560
$ws.Range("H136").Value = 4011.5178
$ws.Range("I136").Value = 2154.7778
$ws.Range("J136").Value = 7353.65
$ws.Range("K136").Value = 6464.3334
$ws.Range("L136").Value = 22060.95
$ws.Range("M136").Value = -3914.3334
$ws.Range("N136").Value = -27160.95
$ws.Range("H137").Value = 55760
$ws.Range("J137").Value = 55760
$ws.Range("L137").Value = 55760
$ws.Range("N137").Value = -65960

$ws = $wb.Worksheets("WVR")
$ws.Range("H20").Value = 10000000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H126").Value = 1216.3334
$ws.Range("I126").Value = 933.6667
$ws.Range("J126").Value = 1499
$ws.Range("K126").Value = 2801.0001
$ws.Range("L126").Value = 4497
$ws.Range("M126").Value = -331.0001000000002
$ws.Range("N126").Value = -9437
$ws.Range("H132").Value = 1220.4117
$ws.Range("I132").Value = 864.95
$ws.Range("J132").Value = 2513
$ws.Range("K132").Value = 2594.85
$ws.Range("L132").Value = 7539
$ws.Range("M132").Value = -64.85000000000036
$ws.Range("N132").Value = -12599
$ws.Range("H133").Value = 40715
$ws.Range("J133").Value = 40715
$ws.Range("L133").Value = 40715
$ws.Range("N133").Value = -50835
$ws.Range("H135").Value = 46471.668
$ws.Range("J135").Value = 46471.668
$ws.Range("L135").Value = 46471.668
$ws.Range("N135").Value = -56611.668
$ws.Range("H136").Value = 3334273.2
$ws.Range("I136").Value = 502.9355
$ws.Range("J136").Value = 8773583
$ws.Range("K136").Value = 1508.8065
$ws.Range("L136").Value = 26320749
$ws.Range("M136").Value = 1041.1935
$ws.Range("N136").Value = -26325849
$ws.Range("H140").Value = 32294
$ws.Range("J140").Value = 32294
$ws.Range("L140").Value = 32294
$ws.Range("N140").Value = -42654

